$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 (Mo., 07.05.2018): fill in begin/end time and break, add a remark
$ws.Range("B46").Value = 0.40625
$ws.Range("C46").Value = 0.70833333333333337
$ws.Range("D46").Value = 1.5
$ws.Range("G46").Value = "Whitebox für Level 1;`nLevel Design mit Artjom besprochen"

# Update the view: scroll position and active selection
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("G47").Select()
